# Auto-generated Excel COM-interop script
# Applies updated market-price figures to the Diabolos_Profits leve-profit sheets
# (columns H,I,J,K,L,M,N) as refreshed by the scheduled data-pull runner.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 325572.94
$ws.Range("J17").Value = 335092.03
$ws.Range("L17").Value = 1005276.09
$ws.Range("N17").Value = -1005612.09
$ws.Range("H28").Value = 43149
$ws.Range("J28").Value = 2392
$ws.Range("L28").Value = 2392
$ws.Range("N28").Value = -3362
$ws.Range("H32").Value = 2431.4443
$ws.Range("J32").Value = 2490.5715
$ws.Range("L32").Value = 2490.5715
$ws.Range("N32").Value = -3142.5715
$ws.Range("H69").Value = 7499.25
$ws.Range("J69").Value = 7499
$ws.Range("L69").Value = 22497
$ws.Range("N69").Value = -24245
$ws.Range("H70").Value = 682765
$ws.Range("I70").Value = 1134730.5
$ws.Range("J70").Value = 4816.6665
$ws.Range("K70").Value = 3404191.5
$ws.Range("L70").Value = 14449.9995
$ws.Range("M70").Value = -3403921.5
$ws.Range("N70").Value = -14989.9995
$ws.Range("H72").Value = 7499.25
$ws.Range("J72").Value = 7499
$ws.Range("L72").Value = 67491
$ws.Range("N72").Value = -76227
$ws.Range("H73").Value = 682765
$ws.Range("I73").Value = 1134730.5
$ws.Range("J73").Value = 4816.6665
$ws.Range("K73").Value = 3404191.5
$ws.Range("L73").Value = 14449.9995
$ws.Range("M73").Value = -3403255.5
$ws.Range("N73").Value = -16321.9995
$ws.Range("H100").Value = 1980.6923
$ws.Range("I100").Value = 1936.9166
$ws.Range("K100").Value = 1936.9166
$ws.Range("M100").Value = -1395.9166
$ws.Range("H116").Value = 97408890
$ws.Range("I116").Value = 62777500
$ws.Range("J116").Value = 166671660
$ws.Range("K116").Value = 62777500
$ws.Range("L116").Value = 166671660
$ws.Range("M116").Value = -62774058
$ws.Range("N116").Value = -166678544
$ws.Range("H132").Value = 29417086
$ws.Range("I132").Value = 31254762
$ws.Range("K132").Value = 93764286
$ws.Range("M132").Value = -93761756
$ws.Range("H138").Value = 3586.7646
$ws.Range("I138").Value = 3196.7646
$ws.Range("J138").Value = 3976.7646
$ws.Range("K138").Value = 9590.293799999999
$ws.Range("L138").Value = 11930.2938
$ws.Range("M138").Value = -4450.293799999999
$ws.Range("N138").Value = -22210.2938

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 21236.492
$ws.Range("I32").Value = 12834.674
$ws.Range("J32").Value = 41577.74
$ws.Range("K32").Value = 12834.674
$ws.Range("L32").Value = 41577.74
$ws.Range("M32").Value = -12547.674
$ws.Range("N32").Value = -42151.74
$ws.Range("H45").Value = 310322.56
$ws.Range("I45").Value = 372069.47
$ws.Range("K45").Value = 372069.47
$ws.Range("M45").Value = -371692.47

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 10190.77
$ws.Range("I20").Value = 2688.8572
$ws.Range("J20").Value = 18943
$ws.Range("K20").Value = 2688.8572
$ws.Range("L20").Value = 18943
$ws.Range("M20").Value = -2441.8572
$ws.Range("N20").Value = -19437
$ws.Range("H99").Value = 3508.8
$ws.Range("I99").Value = 3223.7144
$ws.Range("J99").Value = 7500
$ws.Range("K99").Value = 3223.7144
$ws.Range("L99").Value = 7500
$ws.Range("M99").Value = -1725.7144
$ws.Range("N99").Value = -10496

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 127.63636
$ws.Range("I7").Value = 122.77778
$ws.Range("J7").Value = 149.5
$ws.Range("K7").Value = 122.77778
$ws.Range("L7").Value = 149.5
$ws.Range("M7").Value = -9.777780000000007
$ws.Range("N7").Value = -375.5
$ws.Range("H22").Value = 839.35297
$ws.Range("I22").Value = 418.77777
$ws.Range("J22").Value = 1312.5
$ws.Range("K22").Value = 418.77777
$ws.Range("L22").Value = 1312.5
$ws.Range("M22").Value = -68.77776999999998
$ws.Range("N22").Value = -2012.5
$ws.Range("H31").Value = 3965.0857
$ws.Range("J31").Value = 4587.6816
$ws.Range("L31").Value = 4587.6816
$ws.Range("N31").Value = -5177.6816
$ws.Range("H34").Value = 3965.0857
$ws.Range("J34").Value = 4587.6816
$ws.Range("L34").Value = 4587.6816
$ws.Range("N34").Value = -4991.6816
$ws.Range("H62").Value = 191375.62
$ws.Range("I62").Value = 4701
$ws.Range("K62").Value = 4701
$ws.Range("M62").Value = -4077
$ws.Range("H65").Value = 191375.62
$ws.Range("I65").Value = 4701
$ws.Range("K65").Value = 23505
$ws.Range("M65").Value = -20385
$ws.Range("H99").Value = 405003.84
$ws.Range("J99").Value = 9939.799999999999
$ws.Range("L99").Value = 9939.799999999999
$ws.Range("N99").Value = -12935.8
$ws.Range("H105").Value = 1380.9412
$ws.Range("I105").Value = 1370.5714
$ws.Range("J105").Value = 1429.3334
$ws.Range("K105").Value = 1370.5714
$ws.Range("L105").Value = 1429.3334
$ws.Range("M105").Value = 376.4286
$ws.Range("N105").Value = -4923.3334
$ws.Range("H126").Value = 405003.84
$ws.Range("J126").Value = 9939.799999999999
$ws.Range("L126").Value = 29819.4
$ws.Range("N126").Value = -34759.39999999999
$ws.Range("H132").Value = 219139.9
$ws.Range("I132").Value = 1325.6857
$ws.Range("K132").Value = 3977.0571
$ws.Range("M132").Value = -1447.0571

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 2346.077
$ws.Range("J34").Value = 2458.25
$ws.Range("L34").Value = 7374.75
$ws.Range("N34").Value = -7542.75
$ws.Range("H39").Value = 4049.6428
$ws.Range("J39").Value = 4049.6428
$ws.Range("L39").Value = 12148.9284
$ws.Range("N39").Value = -12736.9284
$ws.Range("H107").Value = 338.875
$ws.Range("J107").Value = 422.2
$ws.Range("L107").Value = 1266.6
$ws.Range("N107").Value = -5106.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 514.9286
$ws.Range("J2").Value = 731.1667
$ws.Range("L2").Value = 731.1667
$ws.Range("N2").Value = -957.1667
$ws.Range("H102").Value = 2884.7778
$ws.Range("I102").Value = 1383
$ws.Range("J102").Value = 4086.2
$ws.Range("K102").Value = 1383
$ws.Range("L102").Value = 4086.2
$ws.Range("M102").Value = 239
$ws.Range("N102").Value = -7330.2
$ws.Range("H113").Value = 3587.7646
$ws.Range("I113").Value = 2822
$ws.Range("K113").Value = 2822
$ws.Range("M113").Value = -652
$ws.Range("H126").Value = 7904.091
$ws.Range("I126").Value = 10845.357
$ws.Range("K126").Value = 32536.071
$ws.Range("M126").Value = -30066.071

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3546.8206
$ws.Range("I46").Value = 1883.4
$ws.Range("K46").Value = 1883.4
$ws.Range("M46").Value = -1695.4
$ws.Range("H70").Value = 55332.332
$ws.Range("J70").Value = 56999
$ws.Range("L70").Value = 56999
$ws.Range("N70").Value = -57539
$ws.Range("H73").Value = 55332.332
$ws.Range("J73").Value = 56999
$ws.Range("L73").Value = 56999
$ws.Range("N73").Value = -58871
$ws.Range("H82").Value = 1944.4
$ws.Range("J82").Value = 1913.2858
$ws.Range("L82").Value = 1913.2858
$ws.Range("N82").Value = -2635.2858
$ws.Range("H85").Value = 1944.4
$ws.Range("J85").Value = 1913.2858
$ws.Range("L85").Value = 1913.2858
$ws.Range("N85").Value = -4409.2858

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H69").Value = 32666.666
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 32666.666
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 32666.666
$ws.Range("M69").ClearContents()
$ws.Range("N69").Value = -34164.666
$ws.Range("H72").Value = 32666.666
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 32666.666
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 97999.99800000001
$ws.Range("M72").ClearContents()
$ws.Range("N72").Value = -105487.998
$ws.Range("H107").Value = 337.5
$ws.Range("I107").Value = 336.42856
$ws.Range("K107").Value = 1009.28568
$ws.Range("M107").Value = 910.71432
$ws.Range("H113").Value = 8019.737
$ws.Range("I113").Value = 12560.556
$ws.Range("J113").Value = 3933
$ws.Range("K113").Value = 37681.66800000001
$ws.Range("L113").Value = 11799
$ws.Range("M113").Value = -35511.66800000001
$ws.Range("N113").Value = -16139
$ws.Range("H132").Value = 435993.38
$ws.Range("I132").Value = 671224.2
$ws.Range("K132").Value = 2013672.6
$ws.Range("M132").Value = -2011142.6
$ws.Range("H136").Value = 8169.278
$ws.Range("I136").Value = 11728.5
$ws.Range("J136").Value = 5321.9
$ws.Range("K136").Value = 35185.5
$ws.Range("M136").Value = -32635.5
$ws.Range("N136").Value = -21065.7
